# Applies the "Recommendation System (Python)" -> "Recommendation server (Python)"
# rename across the two textboxes that still carry the old wording.
#
# Slide 1 (sldId 256), shape id=25 "TextBox 24": the text is split across three
#   runs by the original author ("Recommendation server" + " " + "(Python)").
# Slide 2 (sldId 257), shape id=30 "TextBox 29": simple single-run text swap.

$p = $ppt.ActivePresentation

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# ---- Slide 1: "TextBox 24" (shape id 25) ----
$slide1 = $p.Slides.Item(1)
$shape1 = Get-ShapeById $slide1 25
$tr1 = $shape1.TextFrame.TextRange

# Rewrite the whole run first (keeps a single run briefly)...
$tr1.Text = "Recommendation server (Python)"

# ...then re-assign each segment through Characters() so the paragraph ends up
# split into three runs, matching "Recommendation server" / " " / "(Python)".
$run1 = $tr1.Characters(1, 21)
$run1.Text = "Recommendation server"

$run2 = $tr1.Characters(22, 1)
$run2.Text = " "

$run3 = $tr1.Characters(23, 8)
$run3.Text = "(Python)"

# ---- Slide 2: "TextBox 29" (shape id 30) ----
$slide2 = $p.Slides.Item(2)
$shape2 = Get-ShapeById $slide2 30
$tr2 = $shape2.TextFrame.TextRange
$tr2.Text = "Recommendation server (Python)"
